# Auto-generated edit script: updates market-board derived profit figures
# across the 8 crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# reflecting a refreshed data pull from the scheduled runner.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4433.1665
$ws.Range("I86").Value = 4787.75
$ws.Range("J86").Value = 4149.5
$ws.Range("K86").Value = 4787.75
$ws.Range("L86").Value = 4149.5
$ws.Range("M86").Value = -3664.75
$ws.Range("N86").Value = -6395.5
$ws.Range("H89").Value = 4433.1665
$ws.Range("I89").Value = 4787.75
$ws.Range("J89").Value = 4149.5
$ws.Range("K89").Value = 23938.75
$ws.Range("L89").Value = 20747.5
$ws.Range("M89").Value = -18322.75
$ws.Range("N89").Value = -31979.5
$ws.Range("H98").Value = 2941.8572
$ws.Range("I98").Value = 2610
$ws.Range("K98").Value = 2610
$ws.Range("M98").Value = -1112
$ws.Range("H122").Value = 2941.8572
$ws.Range("I122").Value = 2610
$ws.Range("K122").Value = 7830
$ws.Range("M122").Value = -5380
$ws.Range("H137").Value = 1586.5
$ws.Range("J137").Value = 1564.75
$ws.Range("L137").Value = 4694.25
$ws.Range("N137").Value = -9794.25

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13361.375
$ws.Range("I2").Value = 775.6
$ws.Range("J2").Value = 34337.668
$ws.Range("K2").Value = 775.6
$ws.Range("L2").Value = 34337.668
$ws.Range("M2").Value = -662.6
$ws.Range("N2").Value = -34563.668
$ws.Range("H5").Value = 300
$ws.Range("I5").Value = 300
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 300
$ws.Range("L5").ClearContents()
$ws.Range("N5").Value = 0
$ws.Range("M5").Value = -188
$ws.Range("H32").Value = 2773.397
$ws.Range("I32").Value = 2509.2622
$ws.Range("J32").Value = 5075.143
$ws.Range("K32").Value = 2509.2622
$ws.Range("L32").Value = 5075.143
$ws.Range("M32").Value = -2222.2622
$ws.Range("N32").Value = -5649.143
$ws.Range("H45").Value = 1642.3
$ws.Range("I45").Value = 1802.875
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 1802.875
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -1425.875
$ws.Range("N45").Value = -1754
$ws.Range("H116").Value = 13361.375
$ws.Range("I116").Value = 775.6
$ws.Range("J116").Value = 34337.668
$ws.Range("K116").Value = 775.6
$ws.Range("L116").Value = 34337.668
$ws.Range("M116").Value = 1518.4
$ws.Range("N116").Value = -38925.668
$ws.Range("H132").Value = 2012.8077
$ws.Range("I132").Value = 1550.3684
$ws.Range("J132").Value = 3268
$ws.Range("K132").Value = 4651.1052
$ws.Range("L132").Value = 9804
$ws.Range("M132").Value = -2121.1052
$ws.Range("N132").Value = -14864

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13361.375
$ws.Range("I3").Value = 775.6
$ws.Range("J3").Value = 34337.668
$ws.Range("K3").Value = 775.6
$ws.Range("L3").Value = 34337.668
$ws.Range("M3").Value = -661.6
$ws.Range("N3").Value = -34565.668
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 300
$ws.Range("L4").ClearContents()
$ws.Range("N4").Value = 0
$ws.Range("M4").Value = -185
$ws.Range("H118").Value = 11000
$ws.Range("J118").Value = 11000
$ws.Range("L118").Value = 11000
$ws.Range("N118").Value = -14314
$ws.Range("H134").Value = 8274.210999999999
$ws.Range("I134").Value = 1088.3
$ws.Range("J134").Value = 16258.556
$ws.Range("K134").Value = 3264.9
$ws.Range("L134").Value = 48775.66800000001
$ws.Range("M134").Value = -729.8999999999996
$ws.Range("N134").Value = -53845.66800000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 202.33333
$ws.Range("I7").Value = 51
$ws.Range("K7").Value = 51
$ws.Range("M7").Value = 62
$ws.Range("H31").Value = 1355.4286
$ws.Range("I31").Value = 997.8823
$ws.Range("K31").Value = 997.8823
$ws.Range("M31").Value = -702.8823
$ws.Range("H33").Value = 866.6667
$ws.Range("I33").Value = 866.6667
$ws.Range("K33").Value = 866.6667
$ws.Range("M33").Value = -487.6667
$ws.Range("H34").Value = 1355.4286
$ws.Range("I34").Value = 997.8823
$ws.Range("K34").Value = 997.8823
$ws.Range("M34").Value = -795.8823
$ws.Range("H58").Value = 1609.6666
$ws.Range("I58").Value = 1376.3636
$ws.Range("J58").Value = 1866.3
$ws.Range("K58").Value = 1376.3636
$ws.Range("L58").Value = 1866.3
$ws.Range("M58").Value = -1173.3636
$ws.Range("N58").Value = -2272.3
$ws.Range("H134").Value = 2810.8572
$ws.Range("I134").Value = 3538.75
$ws.Range("J134").Value = 1840.3334
$ws.Range("K134").Value = 10616.25
$ws.Range("L134").Value = 5521.0002
$ws.Range("M134").Value = -8081.25
$ws.Range("N134").Value = -10591.0002
$ws.Range("H136").Value = 1609.6666
$ws.Range("I136").Value = 1376.3636
$ws.Range("J136").Value = 1866.3
$ws.Range("K136").Value = 4129.0908
$ws.Range("L136").Value = 5598.9
$ws.Range("M136").Value = -1579.0908
$ws.Range("N136").Value = -10698.9

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1650746
$ws.Range("I4").Value = 898989
$ws.Range("J4").Value = 1758139.9
$ws.Range("K4").Value = 2696967
$ws.Range("L4").Value = 5274419.699999999
$ws.Range("M4").Value = -2696855
$ws.Range("N4").Value = -5274643.699999999
$ws.Range("H17").Value = 201
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 201
$ws.Range("K17").Value = 0
$ws.Range("L17").ClearContents()
$ws.Range("M17").Value = 603
$ws.Range("N17").Value = -941
$ws.Range("H92").Value = 739.3
$ws.Range("I92").Value = 756.1429000000001
$ws.Range("J92").Value = 700
$ws.Range("K92").Value = 2268.4287
$ws.Range("L92").Value = 2100
$ws.Range("M92").Value = -1020.4287
$ws.Range("N92").Value = -4596
$ws.Range("H131").Value = 13159034
$ws.Range("J131").Value = 1232.9117
$ws.Range("L131").Value = 3698.7351
$ws.Range("N131").Value = -13778.7351

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2428.8
$ws.Range("I102").Value = 2681
$ws.Range("J102").Value = 1420
$ws.Range("K102").Value = 2681
$ws.Range("L102").Value = 1420
$ws.Range("M102").Value = -1059
$ws.Range("N102").Value = -4664
$ws.Range("H132").Value = 3656.3
$ws.Range("I132").Value = 3427.5
$ws.Range("K132").Value = 10282.5
$ws.Range("M132").Value = -7752.5
$ws.Range("H133").Value = 42897
$ws.Range("J133").Value = 42897
$ws.Range("L133").Value = 42897
$ws.Range("N133").Value = -53017

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 21264.549
$ws.Range("I132").Value = 1376.1538
$ws.Range("J132").Value = 41948.48
$ws.Range("K132").Value = 4128.4614
$ws.Range("L132").Value = 125845.44
$ws.Range("M132").Value = -1598.4614
$ws.Range("N132").Value = -130905.44
$ws.Range("H136").Value = 4684
$ws.Range("I136").Value = 5946.4287
$ws.Range("K136").Value = 17839.2861
$ws.Range("M136").Value = -15289.2861

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 34978.57
$ws.Range("J138").Value = 34978.57
$ws.Range("L138").Value = 34978.57
$ws.Range("N138").Value = -45258.57
